$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:w14="http://schemas.microsoft.com/office/word/2010/wordml"'

function Set-ParagraphXml($paraIndex, $innerP) {
    $p = $d.Paragraphs($paraIndex)
    $r = $p.Range
    $r.InsertXML($innerP)
}

# Paragraph 2: "<1,000" + " = 1" runs merge into one run; list numbering/style
# removed in favor of an explicit first-line indent.
$xml2 = '<w:p ' + $wNs + ' w14:paraId="378A63F2" w14:textId="685DA11A" w:rsidR="008E6693" w:rsidRPr="008E6693" w:rsidRDefault="008E6693" w:rsidP="008E6693">'
$xml2 = $xml2 + '<w:pPr>'
$xml2 = $xml2 + '<w:spacing w:after="0" w:line="240" w:lineRule="auto"/>'
$xml2 = $xml2 + '<w:ind w:firstLine="360"/>'
$xml2 = $xml2 + '<w:rPr>'
$xml2 = $xml2 + '<w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/>'
$xml2 = $xml2 + '<w:color w:val="000000"/>'
$xml2 = $xml2 + '</w:rPr>'
$xml2 = $xml2 + '</w:pPr>'
$xml2 = $xml2 + '<w:r w:rsidRPr="008E6693">'
$xml2 = $xml2 + '<w:rPr>'
$xml2 = $xml2 + '<w:rFonts w:ascii="Calibri" w:eastAsia="Times New Roman" w:hAnsi="Calibri" w:cs="Calibri"/>'
$xml2 = $xml2 + '<w:color w:val="000000"/>'
$xml2 = $xml2 + '</w:rPr>'
$xml2 = $xml2 + '<w:t>&lt;1,000 = 1</w:t>'
$xml2 = $xml2 + '</w:r>'
$xml2 = $xml2 + '</w:p>'
Set-ParagraphXml 2 $xml2

# Paragraph 3: "1,001-50,000 = 2" - drop list numbering/style, add first-line indent.
$xml3 = '<w:p ' + $wNs + ' w14:paraId="748ED7A3" w14:textId="21AD59C9" w:rsidR="008E6693" w:rsidRDefault="008E6693" w:rsidP="008E6693">'
$xml3 = $xml3 + '<w:pPr>'
$xml3 = $xml3 + '<w:spacing w:after="0" w:line="240" w:lineRule="auto"/>'
$xml3 = $xml3 + '<w:ind w:firstLine="360"/>'
$xml3 = $xml3 + '</w:pPr>'
$xml3 = $xml3 + '<w:r><w:t>1,001-50,000 = 2</w:t></w:r>'
$xml3 = $xml3 + '</w:p>'
Set-ParagraphXml 3 $xml3

# Paragraph 4: "50,001-250,000 = 3" - same treatment.
$xml4 = '<w:p ' + $wNs + ' w14:paraId="4DEA4B27" w14:textId="479F9BF1" w:rsidR="008E6693" w:rsidRDefault="008E6693" w:rsidP="008E6693">'
$xml4 = $xml4 + '<w:pPr>'
$xml4 = $xml4 + '<w:spacing w:after="0" w:line="240" w:lineRule="auto"/>'
$xml4 = $xml4 + '<w:ind w:firstLine="360"/>'
$xml4 = $xml4 + '</w:pPr>'
$xml4 = $xml4 + '<w:r><w:t>50,001-250,000 = 3</w:t></w:r>'
$xml4 = $xml4 + '</w:p>'
Set-ParagraphXml 4 $xml4

# Paragraph 5: "250,001-1,000,000 = 4" - same treatment, plus the trailing "4"
# is split out into its own run.
$xml5 = '<w:p ' + $wNs + ' w14:paraId="512E506E" w14:textId="60A26FF1" w:rsidR="008E6693" w:rsidRDefault="008E6693" w:rsidP="008E6693">'
$xml5 = $xml5 + '<w:pPr>'
$xml5 = $xml5 + '<w:spacing w:after="0" w:line="240" w:lineRule="auto"/>'
$xml5 = $xml5 + '<w:ind w:firstLine="360"/>'
$xml5 = $xml5 + '</w:pPr>'
$xml5 = $xml5 + '<w:r><w:t xml:space="preserve">250,001-1,000,000 = </w:t></w:r>'
$xml5 = $xml5 + '<w:r><w:t>4</w:t></w:r>'
$xml5 = $xml5 + '</w:p>'
Set-ParagraphXml 5 $xml5

# Paragraph 6: ">1,000,000 = 5" - same treatment, but spacing keeps only
# line/lineRule (no w:after="0").
$xml6 = '<w:p ' + $wNs + ' w14:paraId="1A6DEAC7" w14:textId="17229ECF" w:rsidR="00F15120" w:rsidRDefault="008E6693" w:rsidP="00F15120">'
$xml6 = $xml6 + '<w:pPr>'
$xml6 = $xml6 + '<w:spacing w:line="240" w:lineRule="auto"/>'
$xml6 = $xml6 + '<w:ind w:firstLine="360"/>'
$xml6 = $xml6 + '</w:pPr>'
$xml6 = $xml6 + '<w:r><w:t>&gt;1,000,000 = 5</w:t></w:r>'
$xml6 = $xml6 + '</w:p>'
Set-ParagraphXml 6 $xml6

# Paragraph 12: "Strict Liability" + " = 1" runs merge into a single run;
# numbering and style stay untouched.
$xml12 = '<w:p ' + $wNs + ' w14:paraId="7CDFE7D9" w14:textId="4CA9D40E" w:rsidR="00F15120" w:rsidRDefault="00F15120" w:rsidP="00F15120">'
$xml12 = $xml12 + '<w:pPr>'
$xml12 = $xml12 + '<w:pStyle w:val="ListParagraph"/>'
$xml12 = $xml12 + '<w:numPr>'
$xml12 = $xml12 + '<w:ilvl w:val="0"/>'
$xml12 = $xml12 + '<w:numId w:val="3"/>'
$xml12 = $xml12 + '</w:numPr>'
$xml12 = $xml12 + '</w:pPr>'
$xml12 = $xml12 + '<w:r><w:t>Strict Liability = 1</w:t></w:r>'
$xml12 = $xml12 + '</w:p>'
Set-ParagraphXml 12 $xml12

Write-Output "done"
